$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6733.3
$ws.Range("I40").Value = 6110.5
$ws.Range("J40").Value = 7148.5
$ws.Range("K40").Value = 6110.5
$ws.Range("L40").Value = 7148.5
$ws.Range("M40").Value = -5935.5
$ws.Range("N40").Value = -7498.5
$ws.Range("H88").Value = 1839.875
$ws.Range("I88").Value = 1973.75
$ws.Range("K88").Value = 1973.75
$ws.Range("M88").Value = -1567.75
$ws.Range("H91").Value = 1839.875
$ws.Range("I91").Value = 1973.75
$ws.Range("K91").Value = 1973.75
$ws.Range("M91").Value = -569.75
$ws.Range("H96").Value = 403.75
$ws.Range("I96").Value = 305
$ws.Range("J96").Value = 700
$ws.Range("K96").Value = 915
$ws.Range("L96").Value = 2100
$ws.Range("M96").Value = 458
$ws.Range("N96").Value = -4846
$ws.Range("H103").Value = 27332.666
$ws.Range("I103").Value = 39999.5
$ws.Range("J103").Value = 1999
$ws.Range("K103").Value = 119998.5
$ws.Range("L103").Value = 5997
$ws.Range("M103").Value = -119412.5
$ws.Range("N103").Value = -7169
$ws.Range("H112").Value = 2749.8462
$ws.Range("I112").Value = 3166.6667
$ws.Range("J112").Value = 2624.8
$ws.Range("K112").Value = 9500.000100000001
$ws.Range("L112").Value = 7874.400000000001
$ws.Range("M112").Value = -8392.000100000001
$ws.Range("N112").Value = -10090.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 31735.5
$ws.Range("I28").Value = 31735.5
$ws.Range("K28").Value = 31735.5
$ws.Range("M28").Value = -31543.5
$ws.Range("H44").Value = 12083.5
$ws.Range("J44").Value = 12441.353
$ws.Range("L44").Value = 12441.353
$ws.Range("N44").Value = -13417.353
$ws.Range("H55").Value = 45333
$ws.Range("I55").Value = 15000
$ws.Range("J55").Value = 60499.5
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 60499.5
$ws.Range("M55").Value = -14685
$ws.Range("N55").Value = -61129.5
$ws.Range("H88").Value = 1734.5
$ws.Range("J88").Value = 969
$ws.Range("L88").Value = 969
$ws.Range("N88").Value = -1781
$ws.Range("H91").Value = 1734.5
$ws.Range("J91").Value = 969
$ws.Range("L91").Value = 969
$ws.Range("N91").Value = -3777
$ws.Range("H99").Value = 31735.5
$ws.Range("I99").Value = 31735.5
$ws.Range("K99").Value = 31735.5
$ws.Range("M99").Value = -28740.5
$ws.Range("H102").Value = 7044.8887
$ws.Range("I102").Value = 4601
$ws.Range("K102").Value = 4601
$ws.Range("M102").Value = -2979

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3289.2856
$ws.Range("I20").Value = 3289.2856
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 3289.2856
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -3042.2856
$ws.Range("N20").ClearContents()
$ws.Range("H86").Value = 4955
$ws.Range("I86").Value = 2949.1667
$ws.Range("J86").Value = 8966.666999999999
$ws.Range("K86").Value = 2949.1667
$ws.Range("L86").Value = 8966.666999999999
$ws.Range("M86").Value = -1826.1667
$ws.Range("N86").Value = -11212.667
$ws.Range("H89").Value = 4955
$ws.Range("I89").Value = 2949.1667
$ws.Range("J89").Value = 8966.666999999999
$ws.Range("K89").Value = 14745.8335
$ws.Range("L89").Value = 44833.335
$ws.Range("M89").Value = -9129.833500000001
$ws.Range("N89").Value = -56065.335
$ws.Range("H140").Value = 75093.336
$ws.Range("J140").Value = 75093.336
$ws.Range("L140").Value = 75093.336
$ws.Range("N140").Value = -85453.336

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5250.3213
$ws.Range("I31").Value = 2300.9
$ws.Range("K31").Value = 2300.9
$ws.Range("M31").Value = -2005.9
$ws.Range("H34").Value = 5250.3213
$ws.Range("I34").Value = 2300.9
$ws.Range("K34").Value = 2300.9
$ws.Range("M34").Value = -2098.9
$ws.Range("H55").Value = 16441
$ws.Range("I55").Value = 13999.5
$ws.Range("K55").Value = 13999.5
$ws.Range("M55").Value = -13684.5
$ws.Range("H99").Value = 2372.8
$ws.Range("J99").Value = 2056.25
$ws.Range("L99").Value = 2056.25
$ws.Range("N99").Value = -5052.25
$ws.Range("H122").Value = 1876.8334
$ws.Range("I122").Value = 1912.2
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 5736.6
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -3286.6
$ws.Range("N122").Value = -10000
$ws.Range("H126").Value = 2372.8
$ws.Range("J126").Value = 2056.25
$ws.Range("L126").Value = 6168.75
$ws.Range("N126").Value = -11108.75
$ws.Range("H134").Value = 3887.4443
$ws.Range("I134").Value = 3887.4443
$ws.Range("K134").Value = 11662.3329
$ws.Range("M134").Value = -9127.332900000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 130.66667
$ws.Range("I98").Value = 180
$ws.Range("K98").Value = 540
$ws.Range("M98").Value = 958
$ws.Range("H124").Value = 1111
$ws.Range("J124").Value = 1111
$ws.Range("L124").Value = 3333
$ws.Range("N124").Value = -13153
$ws.Range("H132").Value = 4784.4287
$ws.Range("I132").Value = 4747.3335
$ws.Range("J132").Value = 4812.25
$ws.Range("K132").Value = 42726.0015
$ws.Range("L132").Value = 43310.25
$ws.Range("M132").Value = -40196.0015
$ws.Range("N132").Value = -48370.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 6672.5
$ws.Range("I10").Value = 8866.666999999999
$ws.Range("K10").Value = 8866.666999999999
$ws.Range("M10").Value = -8697.666999999999
$ws.Range("H14").Value = 504
$ws.Range("I14").Value = 504
$ws.Range("K14").Value = 504
$ws.Range("M14").Value = -336
$ws.Range("H99").Value = 3289.8
$ws.Range("I99").Value = 3289.8
$ws.Range("K99").Value = 3289.8
$ws.Range("M99").Value = -1043.8
$ws.Range("H113").Value = 5472.7144
$ws.Range("I113").Value = 3689
$ws.Range("J113").Value = 7094.273
$ws.Range("K113").Value = 3689
$ws.Range("L113").Value = 7094.273
$ws.Range("M113").Value = -1519
$ws.Range("N113").Value = -11434.273
$ws.Range("H126").Value = 4116.7144
$ws.Range("I126").Value = 4153
$ws.Range("J126").Value = 3899
$ws.Range("K126").Value = 12459
$ws.Range("L126").Value = 11697
$ws.Range("M126").Value = -9989
$ws.Range("N126").Value = -16637

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 946
$ws.Range("I22").Value = 697.25
$ws.Range("J22").Value = 1277.6666
$ws.Range("K22").Value = 697.25
$ws.Range("L22").Value = 1277.6666
$ws.Range("M22").Value = -402.25
$ws.Range("N22").Value = -1867.6666
$ws.Range("H27").Value = 946
$ws.Range("I27").Value = 697.25
$ws.Range("J27").Value = 1277.6666
$ws.Range("K27").Value = 697.25
$ws.Range("L27").Value = 1277.6666
$ws.Range("M27").Value = -590.25
$ws.Range("N27").Value = -1491.6666
$ws.Range("H40").Value = 4061.158
$ws.Range("I40").Value = 4398.364
$ws.Range("J40").Value = 3597.5
$ws.Range("K40").Value = 4398.364
$ws.Range("L40").Value = 3597.5
$ws.Range("M40").Value = -4262.364
$ws.Range("N40").Value = -3869.5
$ws.Range("H46").Value = 6728.643
$ws.Range("I46").Value = 2466.6667
$ws.Range("J46").Value = 7891
$ws.Range("K46").Value = 2466.6667
$ws.Range("L46").Value = 7891
$ws.Range("M46").Value = -2278.6667
$ws.Range("N46").Value = -8267
$ws.Range("H74").Value = 46197
$ws.Range("I74").Value = 46197
$ws.Range("K74").Value = 46197
$ws.Range("M74").Value = -45199
$ws.Range("H77").Value = 46197
$ws.Range("I77").Value = 46197
$ws.Range("K77").Value = 138591
$ws.Range("M77").Value = -133599
$ws.Range("H81").Value = 20181
$ws.Range("J81").Value = 20181
$ws.Range("L81").Value = 20181
$ws.Range("N81").Value = -22177
$ws.Range("H84").Value = 20181
$ws.Range("J84").Value = 20181
$ws.Range("L84").Value = 60543
$ws.Range("N84").Value = -70527
$ws.Range("H93").Value = 789.8333
$ws.Range("I93").Value = 445
$ws.Range("J93").Value = 1134.6666
$ws.Range("K93").Value = 445
$ws.Range("L93").Value = 1134.6666
$ws.Range("M93").Value = 803
$ws.Range("N93").Value = -3630.6666
$ws.Range("H122").Value = 4485.143
$ws.Range("I122").Value = 4677.4
$ws.Range("K122").Value = 14032.2
$ws.Range("M122").Value = -11582.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1938
$ws.Range("J3").Value = 1084
$ws.Range("L3").Value = 1084
$ws.Range("N3").Value = -1312
$ws.Range("H11").Value = 2009.8
$ws.Range("J11").Value = 2009.8
$ws.Range("L11").Value = 2009.8
$ws.Range("N11").Value = -2293.8
$ws.Range("H96").Value = 912.8333
$ws.Range("I96").Value = 1082.3334
$ws.Range("J96").Value = 743.3333
$ws.Range("K96").Value = 1082.3334
$ws.Range("L96").Value = 743.3333
$ws.Range("M96").Value = 290.6666
$ws.Range("N96").Value = -3489.3333
$ws.Range("H100").Value = 2040.8
$ws.Range("I100").Value = 2051
$ws.Range("K100").Value = 4102
$ws.Range("M100").Value = -3561
$ws.Range("H136").Value = 3335.375
$ws.Range("I136").Value = 2316.8
$ws.Range("K136").Value = 6950.400000000001
$ws.Range("M136").Value = -4400.400000000001
